# "Final Compile for Fall 2020 Class"
#
# Numerous tweaks made while creating help videos for students:
#   - Fix the "Enrolment Statistics" sheet-name typo -> "Enrollment Statistics"
#     (and keep the chart that lives on that sheet pointing at the renamed sheet).
#   - Make "Enrollment Statistics" the active / selected sheet instead of
#     "Stock Trend".

$wb = $excel.ActiveWorkbook

$oldName = "Enrolment Statistics"
$newName = "Enrollment Statistics"

$sheetNames = @()
foreach ($s in $wb.Worksheets) { $sheetNames += $s.Name }

if ($sheetNames -contains $oldName) {
    $ws = $wb.Worksheets.Item($oldName)
    $ws.Name = $newName
} else {
    $ws = $wb.Worksheets.Item($newName)
}

# Renaming a sheet does not rewrite formulas elsewhere that reference it by name
# (e.g. chart SERIES formulas), so fix those up by hand. This sheet has a pie
# chart whose single series points at the old sheet name.
foreach ($co in $ws.ChartObjects()) {
    $chart = $co.Chart
    foreach ($series in $chart.SeriesCollection()) {
        $formula = $series.Formula
        if ($formula -like "*'$oldName'*") {
            $series.Formula = $formula -replace [regex]::Escape("'$oldName'"), "'$newName'"
        }
    }
}

# This sheet (now "Enrollment Statistics") becomes the active / selected tab,
# replacing "Stock Trend".
$ws.Activate()
